$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (and B/C just in case) to Text format while we write the
# new values, so strings like "1.000" or "0.7009" are not auto-coerced into
# numbers by Excel. We restore the original General format/style afterwards
# so the on-disk cell styling matches the original workbook.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.097.49'
$ws.Range("E2").Value = '  -2.85%  '

$ws.Range("D3").Value = '1.844.20'
$ws.Range("E3").Value = '  -1.84%  '

$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '0.7009'
$ws.Range("E5").Value = '  -5.53%  '

$ws.Range("D6").Value = '236.74'
$ws.Range("E6").Value = '  -2.47%  '

$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").Value = '0.3018'
$ws.Range("E8").Value = '  -4.49%  '

$ws.Range("D9").Value = '0.07346'
$ws.Range("E9").Value = '  +1.64%  '

$ws.Range("D10").Value = '23.21'
$ws.Range("E10").Value = '  -6.31%  '

$ws.Range("D11").Value = '0.08108'
$ws.Range("E11").Value = '  -2.81%  '

$ws.Range("D12").Value = '0.7222'
$ws.Range("E12").Value = '  -4.00%  '

$ws.Range("D13").Value = '1.826.93'
$ws.Range("E13").Value = '  -3.76%  '

$ws.Range("D14").Value = '5.186'
$ws.Range("E14").Value = '  -3.97%  '

$ws.Range("D15").Value = '88.73'
$ws.Range("E15").Value = '  -3.89%  '

$ws.Range("D16").Value = '29.222.56'
$ws.Range("E16").Value = '  -2.55%  '

$ws.Range("D17").Value = '5.756'
$ws.Range("E17").Value = '  -6.23%  '

$ws.Range("D18").Value = '240.19'
$ws.Range("E18").Value = '  -3.30%  '

$ws.Range("D19").Value = '0.000007634'
$ws.Range("E19").Value = '  -2.88%  '

$ws.Range("D20").Value = '12.97'
$ws.Range("E20").Value = '  -4.51%  '

$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.01%  '

$ws.Range("D22").Value = '2.109.76'
$ws.Range("E22").Value = '  -1.15%  '

$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.14%  '

$ws.Range("D24").Value = '7.569'
$ws.Range("E24").Value = '  -5.47%  '

$ws.Range("D25").Value = '0.1467'
$ws.Range("E25").Value = '  -5.71%  '

$ws.Range("D26").Value = '161.73'
$ws.Range("E26").Value = '  -2.72%  '

$ws.Range("D27").Value = '8.956'
$ws.Range("E27").Value = '  -3.74%  '

$ws.Range("D28").Value = '17.99'
$ws.Range("E28").Value = '  -3.71%  '

$ws.Range("D29").Value = '1.923'
$ws.Range("E29").Value = '  -5.42%  '

$ws.Range("E30").Value = '  -7.63%  '

$ws.Range("D31").Value = '4.431'
$ws.Range("E31").Value = '  -3.55%  '

$ws.Range("D32").Value = '1.488'
$ws.Range("E32").Value = '  -3.24%  '

$ws.Range("D33").Value = '3.995'
$ws.Range("E33").Value = '  -5.30%  '

$ws.Range("D34").Value = '0.05181'
$ws.Range("E34").Value = '  -3.47%  '

$ws.Range("E35").Value = '  -5.39%  '

$ws.Range("D36").Value = '0.7086'
$ws.Range("E36").Value = '  -5.94%  '

$ws.Range("D37").Value = '0.9998'
$ws.Range("E37").Value = '  -0.66%  '

$ws.Range("E38").Value = '  -2.24%  '

$ws.Range("D39").Value = '0.01866'
$ws.Range("E39").Value = '  -5.15%  '

$ws.Range("D40").Value = '2.673'
$ws.Range("E40").Value = '  -3.11%  '

$ws.Range("D41").Value = '0.9068'
$ws.Range("E41").Value = '  +5.17%  '

$ws.Range("D42").Value = '0.4267'
$ws.Range("E42").Value = '  -6.18%  '

$ws.Range("D43").Value = '5.891'
$ws.Range("E43").Value = '  -4.31%  '

$ws.Range("D44").Value = '1.055.47'
$ws.Range("E44").Value = '  -6.10%  '

$ws.Range("D45").Value = '69.60'
$ws.Range("E45").Value = '  -4.28%  '

$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  -0.14%  '

$ws.Range("D47").Value = '101.43'
$ws.Range("E47").Value = '  -3.20%  '

$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '1.743'
$ws.Range("E48").Value = '  -6.75%  '

$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '1.997.15'
$ws.Range("E49").Value = '  -1.71%  '

$ws.Range("B50").Value = 'Aptos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D50").Value = '7.073'
$ws.Range("E50").Value = '  -7.17%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '9.157'
$ws.Range("E51").Value = '  -4.05%  '

# Restore original (General / default Normal style) formatting.
$ws.Range("B2:E51").NumberFormat = "General"
$ws.Range("B2:E51").Style = "Normal"